$p = $ppt.ActivePresentation

# --- 1) Slide 1, TextBox 3 (shape id 4): "boxandpointers" -> "__________" ---
$s1 = $p.Slides.Item(1)
$shape = $s1.Shapes.Item("TextBox 3")
$tr = $shape.TextFrame.TextRange
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "__________"

# --- 2) Slide Layout 13 ("1_Title Slide") Date Placeholder: restore old cached date text ---
$master = $p.SlideMaster
$layout13 = $master.CustomLayouts.Item(13)
$dateShape13 = $layout13.Shapes.Item("Date Placeholder 3")
$dateShape13.TextFrame.TextRange.Text = "Tuesday, March 14, 2023"

# --- 3) Slide Layout 14 ("Title and Content") Date Placeholder: restore old cached date text ---
$layout14 = $master.CustomLayouts.Item(14)
$dateShape14 = $layout14.Shapes.Item("Date Placeholder 3")
$dateShape14.TextFrame.TextRange.Text = "Tuesday, March 14, 2023"
